$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("A25").Value = 111936854
$ws.Range("B25").Value = 56414
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 100049
$ws.Range("F25").Value = "Spillkråka"
$ws.Range("G25").Value = "Dryocopus martius"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
$ws.Range("Q25").Value = 450998
$ws.Range("R25").Value = 7087289
$ws.Range("Z25").ClearContents()
$ws.Range("AB25").ClearContents()

# Row 26
$ws.Range("A26").Value = 111936860
$ws.Range("B26").Value = 89423
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 5432
$ws.Range("F26").Value = "Granticka"
$ws.Range("G26").Value = "Porodaedalea chrysoloma"
$ws.Range("H26").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q26").Value = 450975
$ws.Range("R26").Value = 7086983
$ws.Range("Z26").ClearContents()
$ws.Range("AB26").ClearContents()

# Row 27
$ws.Range("A27").Value = 111936894
$ws.Range("B27").Value = 89965
$ws.Range("D27").Value = "VU"
$ws.Range("E27").Value = 760
$ws.Range("F27").Value = "Doftticka"
$ws.Range("G27").Value = "Haploporus odorus"
$ws.Range("H27").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("Q27").Value = 451169
$ws.Range("R27").Value = 7086617
$ws.Range("Z27").ClearContents()
$ws.Range("AB27").ClearContents()

# Row 28
$ws.Range("A28").Value = 111936789
$ws.Range("B28").Value = 90087
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 3298
$ws.Range("F28").Value = "Trådticka"
$ws.Range("G28").Value = "Climacocystis borealis"
$ws.Range("H28").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q28").Value = 450955
$ws.Range("R28").Value = 7087064
$ws.Range("Z28").ClearContents()
$ws.Range("AB28").ClearContents()

# Row 29
$ws.Range("A29").Value = 111936793
$ws.Range("B29").Value = 56398
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = "Tretåig hackspett"
$ws.Range("G29").Value = "Picoides tridactylus"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("Q29").Value = 451089
$ws.Range("R29").Value = 7087233
$ws.Range("Z29").ClearContents()
$ws.Range("AB29").ClearContents()

# Row 30
$ws.Range("A30").Value = 111936892
$ws.Range("B30").Value = 77515
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("Q30").Value = 451172
$ws.Range("R30").Value = 7086727
$ws.Range("Z30").ClearContents()
$ws.Range("AB30").ClearContents()

# Row 31
$ws.Range("A31").Value = 111936864
$ws.Range("B31").Value = 89423
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 5432
$ws.Range("F31").Value = "Granticka"
$ws.Range("G31").Value = "Porodaedalea chrysoloma"
$ws.Range("H31").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q31").Value = 451094
$ws.Range("R31").Value = 7087213
$ws.Range("Z31").ClearContents()
$ws.Range("AB31").ClearContents()

# Move K/L/M/N (empty marker cells) and AC (comment) from row 30 to row 29
$ws.Range("K30:N30").Copy($ws.Range("K29"))
$ws.Range("AC29").Value = $ws.Range("AC30").Value2
$ws.Range("K30:N30").ClearContents()
$ws.Range("AC30").ClearContents()
